# Update the "想去人数" (wanted-to-go count) figures for a handful of events
# on both the "展览" sheet and the aggregated "全部类型" sheet, reflecting the
# refreshed counts from the latest data pull (commit 456a3b4).

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 283
    $ws.Range("F4").Value = 7786
    $ws.Range("F5").Value = 5673
    $ws.Range("F10").Value = 258

    if ($sheetName -eq "展览") {
        $ws.Range("F11").Value = 255
    }
    elseif ($sheetName -eq "全部类型") {
        $ws.Range("F13").Value = 255
    }
}
